$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B5: hours worked on the first day goes from 1 to 2
$ws.Range("B5").Value = 2

# Add a second day row: date 10/23/2012 (serial 41205) with 2 hours
$ws.Range("A6").Value = 41205
$ws.Range("B6").Value = 2

# Move the active selection to B6 to match the post-edit state
$ws.Range("B6").Select()
